$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 172-173, pushing the existing rows
# (old 172-189) down to become rows 174-191. This also grows the
# sheet dimension from R189 to R191 automatically.
$ws.Range("A172:R173").Insert()

# Row 172: new Alcachofa / Espanola record
$ws.Cells.Item(172, 1).Value  = 10
$ws.Cells.Item(172, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(172, 3).Value  = "La Araucanía"
$ws.Cells.Item(172, 4).Value  = 44769
$ws.Cells.Item(172, 5).Value  = 9
$ws.Cells.Item(172, 6).Value  = 100112013
$ws.Cells.Item(172, 7).Value  = "Alcachofa"
$ws.Cells.Item(172, 8).Value  = "Española"
$ws.Cells.Item(172, 9).Value  = "Primera"
$ws.Cells.Item(172, 10).Value = 120
$ws.Cells.Item(172, 11).Value = 18000
$ws.Cells.Item(172, 12).Value = 18000
$ws.Cells.Item(172, 13).Value = 18000
$ws.Cells.Item(172, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(172, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(172, 16).Value = 600
$ws.Cells.Item(172, 17).Value = 30
$ws.Cells.Item(172, 18).Value = "Hortaliza"

# Row 173: new Alcachofa / Madrigal record
$ws.Cells.Item(173, 1).Value  = 10
$ws.Cells.Item(173, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(173, 3).Value  = "La Araucanía"
$ws.Cells.Item(173, 4).Value  = 44769
$ws.Cells.Item(173, 5).Value  = 9
$ws.Cells.Item(173, 6).Value  = 100112013
$ws.Cells.Item(173, 7).Value  = "Alcachofa"
$ws.Cells.Item(173, 8).Value  = "Madrigal"
$ws.Cells.Item(173, 9).Value  = "Primera"
$ws.Cells.Item(173, 10).Value = 200
$ws.Cells.Item(173, 11).Value = 15000
$ws.Cells.Item(173, 12).Value = 15000
$ws.Cells.Item(173, 13).Value = 15000
$ws.Cells.Item(173, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(173, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(173, 16).Value = 375
$ws.Cells.Item(173, 17).Value = 40
$ws.Cells.Item(173, 18).Value = "Hortaliza"
